$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.749
$ws.Range("D4").Value = -7.681999999999999

$ws.Range("D5").Value = -8.206999999999999

$ws.Range("B7").Value = 6.956999999999999

$ws.Range("D8").Value = -7.896000000000001

$ws.Range("B16").Value = 6.427
$ws.Range("D16").Value = -7.917
